$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: frelon_bof ---
$ws.Range("A7").Value = "frelon_bof"
$ws.Range("B7").Value = 161.499
$ws.Range("C7").Value = 253.015
$ws.Range("D7").Value = 398.364
$ws.Range("E7").Value = 479.114
$ws.Range("F7").Value = 559.863
$ws.Range("G7").Value = 640.613
$ws.Range("D7:G7").Interior.ColorIndex = -4142

# --- Row 8: attaque_fre ---
$ws.Range("A8").Value = "attaque_fre"
$ws.Range("B8").Value = 236.865
$ws.Range("C8").Value = 344.531
$ws.Range("D8").Value = 468.347
$ws.Range("B8").Interior.ColorIndex = 6
$ws.Range("D8").Interior.ColorIndex = -4142

# --- Row 9: frelon_nid ---
$ws.Range("A9").Value = "frelon_nid"
$ws.Range("B9").Value = 258.398
$ws.Range("C9").Value = 360.681
$ws.Range("D9").Value = 473.73
$ws.Range("D9").Interior.ColorIndex = 6

# --- Row 10: frelon_nid2 ---
$ws.Range("A10").Value = "frelon_nid2"
$ws.Range("B10").Value = 258.398
$ws.Range("C10").Value = 387.598
$ws.Range("E10").Value = 635.229
$ws.Range("E10").Interior.ColorIndex = -4142

# --- Row 11: attaquev2 ---
$ws.Range("A11").Value = "attaquev2"
$ws.Range("B11").Value = 231.482
$ws.Range("C11").Value = 360.681
$ws.Range("D11").Value = 457.581
$ws.Range("E11").Value = 602.93
$ws.Range("E11").Interior.ColorIndex = -4142

# Match the final selection left by the author
$ws.Range("H11").Select() | Out-Null
